$wb = $excel.ActiveWorkbook

# Overview sheet: Status -> "Ready for handoff", Latest Handoff Date -> "2016-28-17 20:28:15"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-28-17 20:28:15"

# zh-cn sheet: Status -> "Ready for handoff", Latest Handoff Datetime -> "2016-03-17 20:28:12"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-17 20:28:12"

# de-de sheet: Status -> "Ready for handoff", Latest Handoff Datetime -> "2016-03-17 20:28:15"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-17 20:28:15"
